# Refresh market-data columns (H:N) on the leve-profit sheets.
# Source data comes from the Universalis scheduled price-update job;
# this snippet writes the newly-fetched current prices and recomputed
# leve profit figures for the affected rows on each job class sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Cells.Item(86, 8).Value = 27000.75
$ws.Cells.Item(86, 9).Value = 34667.668
$ws.Cells.Item(86, 10).Value = 4000
$ws.Cells.Item(86, 11).Value = 34667.668
$ws.Cells.Item(86, 12).Value = 4000
$ws.Cells.Item(86, 13).Value = -33544.668
$ws.Cells.Item(86, 14).Value = -6246
# Row 89
$ws.Cells.Item(89, 8).Value = 27000.75
$ws.Cells.Item(89, 9).Value = 34667.668
$ws.Cells.Item(89, 10).Value = 4000
$ws.Cells.Item(89, 11).Value = 173338.34
$ws.Cells.Item(89, 12).Value = 20000
$ws.Cells.Item(89, 13).Value = -167722.34
$ws.Cells.Item(89, 14).Value = -31232
# Row 98
$ws.Cells.Item(98, 8).Value = 2719.0588
$ws.Cells.Item(98, 9).Value = 2403.2
$ws.Cells.Item(98, 10).Value = 5088
$ws.Cells.Item(98, 11).Value = 2403.2
$ws.Cells.Item(98, 12).Value = 5088
$ws.Cells.Item(98, 13).Value = -905.1999999999998
$ws.Cells.Item(98, 14).Value = -8084
# Row 101
$ws.Cells.Item(101, 8).Value = 1440.625
$ws.Cells.Item(101, 9).Value = 591.6667
$ws.Cells.Item(101, 11).Value = 1775.0001
$ws.Cells.Item(101, 13).Value = -153.0001
# Row 116
$ws.Cells.Item(116, 8).Value = 3784.3547
$ws.Cells.Item(116, 9).Value = 2764.9412
$ws.Cells.Item(116, 10).Value = 5022.2144
$ws.Cells.Item(116, 11).Value = 2764.9412
$ws.Cells.Item(116, 12).Value = 5022.2144
$ws.Cells.Item(116, 13).Value = 677.0587999999998
$ws.Cells.Item(116, 14).Value = -11906.2144
# Row 122
$ws.Cells.Item(122, 8).Value = 2719.0588
$ws.Cells.Item(122, 9).Value = 2403.2
$ws.Cells.Item(122, 10).Value = 5088
$ws.Cells.Item(122, 11).Value = 7209.599999999999
$ws.Cells.Item(122, 12).Value = 15264
$ws.Cells.Item(122, 13).Value = -4759.599999999999
$ws.Cells.Item(122, 14).Value = -20164
# Row 137
$ws.Cells.Item(137, 8).Value = 5004789
$ws.Cells.Item(137, 9).Value = 8339449
$ws.Cells.Item(137, 10).Value = 2799.875
$ws.Cells.Item(137, 11).Value = 25018347
$ws.Cells.Item(137, 12).Value = 8399.625
$ws.Cells.Item(137, 13).Value = -25015797
$ws.Cells.Item(137, 14).Value = -13499.625
# Row 138
$ws.Cells.Item(138, 8).Value = 2690.7058
$ws.Cells.Item(138, 9).Value = 1370.2593
$ws.Cells.Item(138, 10).Value = 4176.2085
$ws.Cells.Item(138, 11).Value = 4110.7779
$ws.Cells.Item(138, 12).Value = 12528.6255
$ws.Cells.Item(138, 13).Value = 1029.2221
$ws.Cells.Item(138, 14).Value = -22808.6255

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5900.8716
$ws.Cells.Item(32, 9).Value = 4660.396
$ws.Cells.Item(32, 10).Value = 9768.235000000001
$ws.Cells.Item(32, 11).Value = 4660.396
$ws.Cells.Item(32, 12).Value = 9768.235000000001
$ws.Cells.Item(32, 13).Value = -4373.396
$ws.Cells.Item(32, 14).Value = -10342.235
# Row 61
$ws.Cells.Item(61, 8).Value = 2251.3142
$ws.Cells.Item(61, 9).Value = 1371.3
$ws.Cells.Item(61, 10).Value = 3424.6667
$ws.Cells.Item(61, 11).Value = 1371.3
$ws.Cells.Item(61, 12).Value = 3424.6667
$ws.Cells.Item(61, 13).Value = -1159.3
$ws.Cells.Item(61, 14).Value = -3848.6667
# Row 74
$ws.Cells.Item(74, 8).Value = 604.6667
$ws.Cells.Item(74, 9).Value = 572.1667
$ws.Cells.Item(74, 10).Value = 799.6667
$ws.Cells.Item(74, 11).Value = 572.1667
$ws.Cells.Item(74, 12).Value = 799.6667
$ws.Cells.Item(74, 13).Value = 301.8333
$ws.Cells.Item(74, 14).Value = -2547.6667
# Row 77
$ws.Cells.Item(77, 8).Value = 604.6667
$ws.Cells.Item(77, 9).Value = 572.1667
$ws.Cells.Item(77, 10).Value = 799.6667
$ws.Cells.Item(77, 11).Value = 2860.8335
$ws.Cells.Item(77, 12).Value = 3998.3335
$ws.Cells.Item(77, 13).Value = 1507.1665
$ws.Cells.Item(77, 14).Value = -12734.3335
# Row 132
$ws.Cells.Item(132, 8).Value = 22729916
$ws.Cells.Item(132, 9).Value = 33335332
$ws.Cells.Item(132, 11).Value = 100005996
$ws.Cells.Item(132, 13).Value = -100003466
# Row 136
$ws.Cells.Item(136, 8).Value = 2251.3142
$ws.Cells.Item(136, 9).Value = 1371.3
$ws.Cells.Item(136, 10).Value = 3424.6667
$ws.Cells.Item(136, 11).Value = 4113.9
$ws.Cells.Item(136, 12).Value = 10274.0001
$ws.Cells.Item(136, 13).Value = -1563.9
$ws.Cells.Item(136, 14).Value = -15374.0001

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 2106.4348
$ws.Cells.Item(107, 9).Value = 1499.5834
$ws.Cells.Item(107, 10).Value = 2768.4546
$ws.Cells.Item(107, 11).Value = 1499.5834
$ws.Cells.Item(107, 12).Value = 2768.4546
$ws.Cells.Item(107, 13).Value = 420.4166
$ws.Cells.Item(107, 14).Value = -6608.4546
# Row 134
$ws.Cells.Item(134, 8).Value = 2790.75
$ws.Cells.Item(134, 9).Value = 1988
$ws.Cells.Item(134, 10).Value = 4556.8
$ws.Cells.Item(134, 11).Value = 5964
$ws.Cells.Item(134, 12).Value = 13670.4
$ws.Cells.Item(134, 13).Value = -3429
$ws.Cells.Item(134, 14).Value = -18740.4

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2566982
$ws.Cells.Item(31, 9).Value = 2780174.8
$ws.Cells.Item(31, 10).Value = 8666.666999999999
$ws.Cells.Item(31, 11).Value = 2780174.8
$ws.Cells.Item(31, 12).Value = 8666.666999999999
$ws.Cells.Item(31, 13).Value = -2779879.8
$ws.Cells.Item(31, 14).Value = -9256.666999999999
# Row 34
$ws.Cells.Item(34, 8).Value = 2566982
$ws.Cells.Item(34, 9).Value = 2780174.8
$ws.Cells.Item(34, 10).Value = 8666.666999999999
$ws.Cells.Item(34, 11).Value = 2780174.8
$ws.Cells.Item(34, 12).Value = 8666.666999999999
$ws.Cells.Item(34, 13).Value = -2779972.8
$ws.Cells.Item(34, 14).Value = -9070.666999999999
# Row 58
$ws.Cells.Item(58, 8).Value = 16669593
$ws.Cells.Item(58, 9).Value = 1697.7858
$ws.Cells.Item(58, 10).Value = 31254002
$ws.Cells.Item(58, 11).Value = 1697.7858
$ws.Cells.Item(58, 12).Value = 31254002
$ws.Cells.Item(58, 13).Value = -1494.7858
$ws.Cells.Item(58, 14).Value = -31254408
# Row 70
$ws.Cells.Item(70, 8).Value = 34374.875
$ws.Cells.Item(70, 10).Value = 34285.715
$ws.Cells.Item(70, 12).Value = 34285.715
$ws.Cells.Item(70, 14).Value = -34915.715
# Row 73
$ws.Cells.Item(73, 8).Value = 34374.875
$ws.Cells.Item(73, 10).Value = 34285.715
$ws.Cells.Item(73, 12).Value = 34285.715
$ws.Cells.Item(73, 14).Value = -36469.715
# Row 132
$ws.Cells.Item(132, 8).Value = 3921.15
$ws.Cells.Item(132, 9).Value = 2840.2307
$ws.Cells.Item(132, 11).Value = 8520.6921
$ws.Cells.Item(132, 13).Value = -5990.6921
# Row 134
$ws.Cells.Item(134, 8).Value = 1366.1389
$ws.Cells.Item(134, 9).Value = 818.5
$ws.Cells.Item(134, 10).Value = 2790
$ws.Cells.Item(134, 11).Value = 2455.5
$ws.Cells.Item(134, 12).Value = 8370
$ws.Cells.Item(134, 13).Value = 79.5
$ws.Cells.Item(134, 14).Value = -13440
# Row 136
$ws.Cells.Item(136, 8).Value = 16669593
$ws.Cells.Item(136, 9).Value = 1697.7858
$ws.Cells.Item(136, 10).Value = 31254002
$ws.Cells.Item(136, 11).Value = 5093.357400000001
$ws.Cells.Item(136, 12).Value = 93762006
$ws.Cells.Item(136, 13).Value = -2543.357400000001
$ws.Cells.Item(136, 14).Value = -93767106

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Cells.Item(80, 8).Value = 4117.1665
$ws.Cells.Item(80, 9).Value = 1400
$ws.Cells.Item(80, 10).Value = 5475.75
$ws.Cells.Item(80, 11).Value = 4200
$ws.Cells.Item(80, 12).Value = 16427.25
$ws.Cells.Item(80, 13).Value = -3264
$ws.Cells.Item(80, 14).Value = -18299.25
# Row 83
$ws.Cells.Item(83, 8).Value = 4117.1665
$ws.Cells.Item(83, 9).Value = 1400
$ws.Cells.Item(83, 10).Value = 5475.75
$ws.Cells.Item(83, 11).Value = 12600
$ws.Cells.Item(83, 12).Value = 49281.75
$ws.Cells.Item(83, 13).Value = -7920
$ws.Cells.Item(83, 14).Value = -58641.75
# Row 117
$ws.Cells.Item(117, 8).Value = 1144.875
$ws.Cells.Item(117, 10).Value = 1222.7142
$ws.Cells.Item(117, 12).Value = 3668.1426
$ws.Cells.Item(117, 14).Value = -10552.1426
# Row 131
$ws.Cells.Item(131, 8).Value = 1040.2632
$ws.Cells.Item(131, 9).Value = 608.3333
$ws.Cells.Item(131, 10).Value = 1121.25
$ws.Cells.Item(131, 11).Value = 1824.9999
$ws.Cells.Item(131, 12).Value = 3363.75
$ws.Cells.Item(131, 13).Value = 3215.0001
$ws.Cells.Item(131, 14).Value = -13443.75
# Row 132
$ws.Cells.Item(132, 8).Value = 4320.8
$ws.Cells.Item(132, 10).Value = 5301
$ws.Cells.Item(132, 12).Value = 47709
$ws.Cells.Item(132, 14).Value = -52769
# Row 133
$ws.Cells.Item(133, 8).Value = 4619.3184
$ws.Cells.Item(133, 9).Value = 7515
$ws.Cells.Item(133, 10).Value = 4329.75
$ws.Cells.Item(133, 11).Value = 22545
$ws.Cells.Item(133, 12).Value = 12989.25
$ws.Cells.Item(133, 13).Value = -17485
$ws.Cells.Item(133, 14).Value = -23109.25

$ws = $wb.Worksheets.Item("GSM")
# Row 69
$ws.Cells.Item(69, 8).Value = 31562.5
$ws.Cells.Item(69, 10).Value = 31562.5
$ws.Cells.Item(69, 12).Value = 31562.5
$ws.Cells.Item(69, 14).Value = -33060.5
# Row 72
$ws.Cells.Item(72, 8).Value = 31562.5
$ws.Cells.Item(72, 10).Value = 31562.5
$ws.Cells.Item(72, 12).Value = 94687.5
$ws.Cells.Item(72, 14).Value = -102175.5
# Row 132
$ws.Cells.Item(132, 8).Value = 3282.2896
$ws.Cells.Item(132, 9).Value = 2857.7827
$ws.Cells.Item(132, 11).Value = 8573.348100000001
$ws.Cells.Item(132, 13).Value = -6043.348100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Cells.Item(122, 8).Value = 2845.6858
$ws.Cells.Item(122, 9).Value = 2475.2
$ws.Cells.Item(122, 10).Value = 3339.6667
$ws.Cells.Item(122, 11).Value = 7425.599999999999
$ws.Cells.Item(122, 12).Value = 10019.0001
$ws.Cells.Item(122, 13).Value = -4975.599999999999
$ws.Cells.Item(122, 14).Value = -14919.0001
# Row 132
$ws.Cells.Item(132, 8).Value = 3225.0715
$ws.Cells.Item(132, 9).Value = 2010.4
$ws.Cells.Item(132, 10).Value = 3899.889
$ws.Cells.Item(132, 11).Value = 6031.200000000001
$ws.Cells.Item(132, 12).Value = 11699.667
$ws.Cells.Item(132, 13).Value = -3501.200000000001
$ws.Cells.Item(132, 14).Value = -16759.667
# Row 136
$ws.Cells.Item(136, 8).Value = 4765309.5
$ws.Cells.Item(136, 9).Value = 11114811
$ws.Cells.Item(136, 10).Value = 3183.6667
$ws.Cells.Item(136, 11).Value = 33344433
$ws.Cells.Item(136, 12).Value = 9551.000100000001
$ws.Cells.Item(136, 13).Value = -33341883
$ws.Cells.Item(136, 14).Value = -14651.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 1734
$ws.Cells.Item(62, 9).Value = 1734
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 1734
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -1110
$ws.Cells.Item(62, 14).ClearContents()
# Row 65
$ws.Cells.Item(65, 8).Value = 1734
$ws.Cells.Item(65, 9).Value = 1734
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 8670
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -5550
$ws.Cells.Item(65, 14).ClearContents()
# Row 126
$ws.Cells.Item(126, 8).Value = 3848146.5
$ws.Cells.Item(126, 9).Value = 1561.0555
$ws.Cells.Item(126, 10).Value = 12502963
$ws.Cells.Item(126, 11).Value = 4683.166499999999
$ws.Cells.Item(126, 12).Value = 37508889
$ws.Cells.Item(126, 13).Value = -2213.166499999999
$ws.Cells.Item(126, 14).Value = -37513829
# Row 132
$ws.Cells.Item(132, 8).Value = 305433.44
$ws.Cells.Item(132, 9).Value = 386242.5
$ws.Cells.Item(132, 11).Value = 1158727.5
$ws.Cells.Item(132, 13).Value = -1156197.5
# Row 136
$ws.Cells.Item(136, 8).Value = 2091.8
$ws.Cells.Item(136, 9).Value = 859.625
$ws.Cells.Item(136, 11).Value = 2578.875
$ws.Cells.Item(136, 13).Value = -28.875

